# Update cryptos list (prices / volume percentages) as scraped on
# Tue Nov 14 22:39:58 UTC 2023 with GitHub Actions.
#
# Price cells (column D) are stored as plain text in the source sheet
# (e.g. "35.543.34", "241.69"). Any value that would otherwise be
# re-interpreted by Excel as a number (single "." separator) is forced
# back to text via NumberFormat "@" so the stored cell keeps its
# original text nature instead of silently becoming numeric.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "35.543.34"
$ws.Range("E2").Value = "  -2.95%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "1.980.77"
$ws.Range("E3").Value = "  -4.23%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  +0.10%  "

# --- Row 5: BNB ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.69"
$ws.Range("E5").Value = "  -0.35%  "

# --- Row 6: XRP ---
$ws.Range("E6").Value = "  -3.83%  "

# --- Row 7: USDC ---
$ws.Range("E7").Value = "  +0.02%  "

# --- Row 8: Solana ---
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.23"
$ws.Range("E8").Value = "  +6.12%  "

# --- Row 9: OKB ---
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.19"
$ws.Range("E9").Value = "  +1.79%  "

# --- Row 10: Cardano ---
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.357"
$ws.Range("E10").Value = "  -0.66%  "

# --- Row 11: Dogecoin ---
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0727"
$ws.Range("E11").Value = "  -3.46%  "

# --- Row 12: TRON ---
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.102"
$ws.Range("E12").Value = "  -5.28%  "

# --- Row 13: Polygon ---
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.908"
$ws.Range("E13").Value = "  -0.12%  "

# --- Row 14: Chainlink ---
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.08"
$ws.Range("E14").Value = "  -4.13%  "

# --- Row 15: WrappedliquidstakedEther2.0 ---
$ws.Range("D15").Value = "2.273.04"
$ws.Range("E15").Value = "  -4.05%  "

# --- Row 16: Polkadot ---
$ws.Range("E16").Value = "  -3.66%  "

# --- Row 17: WrappedEther ---
$ws.Range("D17").Value = "1.984.99"
$ws.Range("E17").Value = "  -4.62%  "

# --- Row 18: Avalanche ---
$ws.Range("E18").Value = "  +3.55%  "

# --- Row 19: WrappedBTC ---
$ws.Range("D19").Value = "35.476.35"
$ws.Range("E19").Value = "  -2.93%  "

# --- Row 20: Litecoin ---
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.09"
$ws.Range("E20").Value = "  -2.47%  "

# --- Row 21: ShibaInu ---
$ws.Range("E21").Value = "  -3.56%  "

# --- Row 22: BitcoinCash ---
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "232.31"
$ws.Range("E22").Value = "  -2.09%  "

# --- Row 23: Uniswap ---
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.05"
$ws.Range("E23").Value = "  -4.50%  "

# --- Row 24: Dai ---
$ws.Range("E24").Value = "  -0.06%  "

# --- Row 25: Toncoin ---
$ws.Range("E25").Value = "  -3.18%  "

# --- Row 26: PancakeSwap ---
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.30"
$ws.Range("E26").Value = "  +7.84%  "

# --- Row 27: Monero ---
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.35"
$ws.Range("E27").Value = "  -0.88%  "

# --- Row 28: Cosmos ---
$ws.Range("E28").Value = "  -4.21%  "

# --- Row 29: EthereumClassic ---
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.41"
$ws.Range("E29").Value = "  -5.67%  "

# --- Row 30: Stellar ---
$ws.Range("E30").Value = "  -3.11%  "

# --- Row 31: Filecoin ---
$ws.Range("E31").Value = "  -5.80%  "

# --- Row 32: ImmutableX ---
$ws.Range("E32").Value = "  -2.36%  "

# --- Row 33: Hedera ---
$ws.Range("E33").Value = "  -2.28%  "

# --- Row 34: Kaspa ---
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0909"
$ws.Range("E34").Value = "  +10.82%  "

# --- Row 35: InternetComputer(DFINITY) ---
$ws.Range("E35").Value = "  -7.91%  "

# --- Row 36: LidoDAOToken ---
$ws.Range("E36").Value = "  +0.79%  "

# --- Row 37: BinanceUSD ---
$ws.Range("E37").Value = "  +0.11%  "

# --- Row 38: WEMIXToken ---
$ws.Range("E38").Value = "  -2.34%  "

# --- Row 39: THORChain ---
$ws.Range("E39").Value = "  +0.01%  "

# --- Row 40: TrustWalletToken ---
$ws.Range("E40").Value = "  -5.74%  "

# --- Row 41: HuobiToken ---
$ws.Range("E41").Value = "  -3.49%  "

# --- Row 42: VeChain ---
$ws.Range("E42").Value = "  -3.47%  "

# --- Row 43: ARBITRUM ---
$ws.Range("E43").Value = "  -5.14%  "

# --- Row 44: Cronos ---
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0886"
$ws.Range("E44").Value = "  -5.22%  "

# --- Rows 45 & 46: Aave and Maker swap ranking positions ---
# Row 45 becomes Maker, Row 46 becomes Aave, with refreshed prices.
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.375.53"
$ws.Range("E45").Value = "  -1.56%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.46"
$ws.Range("E46").Value = "  -4.31%  "

# --- Row 47: FraxShare ---
$ws.Range("E47").Value = "  -1.58%  "

# --- Row 48: InjectiveProtocol ---
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.39"
$ws.Range("E48").Value = "  -1.45%  "

# --- Row 49: MXToken ---
$ws.Range("E49").Value = "  +0.35%  "

# --- Row 50: RenderToken ---
$ws.Range("E50").Value = "  -4.46%  "

# --- Row 51: MultiversX ---
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.60"
$ws.Range("E51").Value = "  +0.55%  "
